# TC07_Canine_StudyUBC02-Breed_Diagnosis_PrimDiseaseSite.xlsx
#
# The "CasesTab" row's Cypher query (cell B2 on the "startup" sheet) is
# updated: the `Cohort` output column (and the trailing comma that used to
# precede it on the `Response to Treatment` line) is removed from the
# RETURN clause. No other cell content changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("startup")

$newCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n`nMATCH (c)<--(diag:diagnosis)`nWHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Mixed Breed', 'Scottish Terrier','Shetland Sheepdog']and diag.disease_term in ['Bladder Cancer','Healthy Control'] and diag.primary_disease_site in ['Bladder', 'Bladder, Urethra', 'Bladder, Urethra, Vagina']`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newCasesQuery

# The rows re-wrap to one fewer line (Cases) / reflow at the new line
# height once re-opened, so the row heights shrink accordingly.
$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(3).RowHeight = 275.5
$ws.Rows.Item(4).RowHeight = 290

# Selection moved up to the (now shorter) Cases-query cell.
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
